$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value2 = 5
$ws.Range("C1").Value2 = 2.5
$ws.Range("D1").Value2 = 1.25
$ws.Range("E1").Value2 = 0.625
$ws.Range("F1").Value2 = 0.3125
$ws.Range("G1").Value2 = 0.15625
$ws.Range("H1").Value2 = 0.078125
$ws.Range("I1").Value2 = 0.0390625
$ws.Range("A2").Value2 = 36
$ws.Range("B2").Value2 = 0.999663589130191
$ws.Range("C2").Value2 = 0.9999626210144656
$ws.Range("D2").Value2 = 0.9999626210144656
$ws.Range("E2").Value2 = 0.9999626210144656
$ws.Range("F2").Value2 = 0.9963742384031697
$ws.Range("G2").Value2 = 0.9999626210144656
$ws.Range("H2").Value2 = 0.5409486786528613
$ws.Range("I2").Value2 = 0.5520128583710238
$ws.Range("J2").Value2 = 0.6216872874070197
$ws.Range("A3").Value2 = 18
$ws.Range("B3").Value2 = 1.0002616528987402
$ws.Range("C3").Value2 = 1.000560684783015
$ws.Range("D3").Value2 = 1.0002616528987402
$ws.Range("E3").Value2 = 1.000560684783015
$ws.Range("F3").Value2 = 0.9984674615930923
$ws.Range("G3").Value2 = 0.43120397712406083
$ws.Range("H3").Value2 = 0.3687063133106568
$ws.Range("I3").Value2 = 0.3755840466489739
$ws.Range("J3").Value2 = 0.6192950323328223
$ws.Range("A4").Value2 = 9
$ws.Range("B4").Value2 = 0.9999626210144656
$ws.Range("C4").Value2 = 0.9999626210144656
$ws.Range("D4").Value2 = 0.999663589130191
$ws.Range("E4").Value2 = 0.9437446267708294
$ws.Range("F4").Value2 = 0.4383807423466527
$ws.Range("G4").Value2 = 0.40877658580346127
$ws.Range("H4").Value2 = 0.34657795387433177
$ws.Range("I4").Value2 = 0.35794116547676896
$ws.Range("J4").Value2 = 0.5699547714275034
$ws.Range("A5").Value2 = 4.5
$ws.Range("B5").Value2 = 0.9999626210144656
$ws.Range("C5").Value2 = 0.998766493477367
$ws.Range("D5").Value2 = 0.995178110866071
$ws.Range("E5").Value2 = 0.6097260120360333
$ws.Range("F5").Value2 = 0.4267184988599409
$ws.Range("G5").Value2 = 0.3974133742010241
$ws.Range("H5").Value2 = 0.31727282921541505
$ws.Range("I5").Value2 = 0.3582401973610436
$ws.Range("J5").Value2 = 0.605539565656188
$ws.Range("A6").Value2 = 2.25
$ws.Range("B6").Value2 = 0.998766493477367
$ws.Range("C6").Value2 = 1.000560684783015
$ws.Range("D6").Value2 = 0.9383620528538855
$ws.Range("E6").Value2 = 0.5436399656113332
$ws.Range("F6").Value2 = 0.4341942959668074
$ws.Range("G6").Value2 = 0.3426905393787612
$ws.Range("H6").Value2 = 0.27331514222703984
$ws.Range("I6").Value2 = 0.3325234553134228
$ws.Range("J6").Value2 = 0.6121182671102305
$ws.Range("A7").Value2 = 1.125
$ws.Range("B7").Value2 = 0.9999626210144656
$ws.Range("C7").Value2 = 0.999663589130191
$ws.Range("D7").Value2 = 0.7048181512353755
$ws.Range("E7").Value2 = 0.47755391918663326
$ws.Range("F7").Value2 = 0.3785743654917206
$ws.Range("G7").Value2 = 0.3558479422868463
$ws.Range("H7").Value2 = 0.3193660524053377
$ws.Range("I7").Value2 = 0.3029192987702313
$ws.Range("J7").Value2 = 0.5532089859081224
$ws.Range("A8").Value2 = 0.5625
$ws.Range("B8").Value2 = 0.9990655253616416
$ws.Range("C8").Value2 = 0.9915897282547751
$ws.Range("D8").Value2 = 0.5056629163084513
$ws.Range("E8").Value2 = 0.456621687287407
$ws.Range("F8").Value2 = 0.2484954958322431
$ws.Range("G8").Value2 = 0.2667364407729974
$ws.Range("H8").Value2 = 0.27062385526856797
$ws.Range("I8").Value2 = 0.25806451612903225
$ws.Range("J8").Value2 = 0.53765932792584
$ws.Range("A9").Value2 = 0.28125
$ws.Range("B9").Value2 = 1.002354876088663
$ws.Range("C9").Value2 = 1.0035510036257616
$ws.Range("D9").Value2 = 1.0038500355100362
$ws.Range("E9").Value2 = 1.0035510036257616
$ws.Range("F9").Value2 = 0.24072066684110185
$ws.Range("G9").Value2 = 0.2876686726722236
$ws.Range("H9").Value2 = 0.2966396292004635
$ws.Range("I9").Value2 = 0.3364108698089934
$ws.Range("J9").Value2 = 0.53765932792584
$ws.Range("B10").Value2 = 1.000560684783015
$ws.Range("C10").Value2 = 0.983814899263634
$ws.Range("D10").Value2 = 0.5576944641722423
$ws.Range("E10").Value2 = 0.44017493365230065
$ws.Range("F10").Value2 = 0.3097970321085486
$ws.Range("G10").Value2 = 0.3657159944679101
$ws.Range("H10").Value2 = 0.3097970321085486
$ws.Range("I10").Value2 = 0.3097970321085486

Write-Host "Done"